# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund-holdings detail) right before the
#   "总计" (summary) sheet.
# - Prepend a "2022-Q1" row to the "总计" summary table (shifting the
#   existing rows down by one).

$wb = $excel.ActiveWorkbook

# A worksheet we already know carries the "normal" look (bold header /
# bordered first column, style index reused rather than re-created) that
# every quarterly sheet in this workbook shares. We use it purely as a
# formatting donor via Copy/PasteSpecial so no duplicate styles get baked
# into styles.xml.
$fmtDonor = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Remember the existing "总计" table so we can rebuild it (shifted).
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")

$oldDates  = @()
$oldCounts = @()
$oldValues = @()
$r = 2
while ($true) {
    $dateVal = $oldTotal.Cells.Item($r, 2).Value2
    if ($dateVal -eq $null) { break }
    $oldDates  += $dateVal
    $oldCounts += $oldTotal.Cells.Item($r, 3).Value2
    $oldValues += $oldTotal.Cells.Item($r, 4).Value2
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Remove the old "总计" sheet, then re-add "2022-Q1" and "总计" in the
#    right order so the sheetId allocation lines up (2022-Q1 -> 6,
#    总计 -> 7) the same way Excel would after a delete + two fresh adds.
# ---------------------------------------------------------------------
$oldTotal.Delete()

$q1 = $wb.Worksheets.Add($null, $fmtDonor)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# 3. Fill in "2022-Q1" (fund holdings detail).
# ---------------------------------------------------------------------
$fmtDonor.Range("B1:H1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"

$fmtDonor.Range("A2:A4").Copy() | Out-Null
$q1.Range("A2:A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$q1Rows = @(
    @(0, "000586", "景顺长城中小板创业板精选股票", "2.42", "94.15", "6.83", "0.1653", 3),
    @(1, "010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "6.22", "0.0703", 4),
    @(2, "260115", "景顺长城中小盘混合", "0.96", "94.00", "5.56", "0.0534", 5)
)

$row = 2
foreach ($rec in $q1Rows) {
    $q1.Cells.Item($row, 1).Value2 = $rec[0]

    # Fund code ("000586", …) must stay text — a leading zero would
    # otherwise be silently dropped by numeric auto-detection.
    $c = $q1.Cells.Item($row, 2); $c.NumberFormat = "@"; $c.Value2 = $rec[1]

    $q1.Cells.Item($row, 3).Value2 = $rec[2]

    $c = $q1.Cells.Item($row, 4); $c.NumberFormat = "@"; $c.Value2 = $rec[3]
    $c = $q1.Cells.Item($row, 5); $c.NumberFormat = "@"; $c.Value2 = $rec[4]
    $c = $q1.Cells.Item($row, 6); $c.NumberFormat = "@"; $c.Value2 = $rec[5]
    $c = $q1.Cells.Item($row, 7); $c.NumberFormat = "@"; $c.Value2 = $rec[6]

    $q1.Cells.Item($row, 8).Value2 = $rec[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Rebuild "总计" with a new "2022-Q1" row on top of the old data.
# ---------------------------------------------------------------------
$fmtDonor.Range("B1:D1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$total.Range("B1").Value2 = "日期"
$total.Range("C1").Value2 = "持有数量(只)"
$total.Range("D1").Value2 = "持有市值(亿元)"

$newDates  = @("2022-Q1") + $oldDates
$newCounts = @(3) + $oldCounts
$newValues = @(0.29) + $oldValues

$lastRow = $newDates.Count + 1
$fmtDonor.Range("A2").Copy() | Out-Null
$total.Range("A2:A$lastRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$row = 2
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $total.Cells.Item($row, 1).Value2 = $i
    $total.Cells.Item($row, 2).Value2 = $newDates[$i]
    $total.Cells.Item($row, 3).Value2 = $newCounts[$i]
    $total.Cells.Item($row, 4).Value2 = $newValues[$i]
    $row = $row + 1
}

$q1.Range("A1").Select() | Out-Null
$total.Range("A1").Select() | Out-Null
